# ---------------------------------------------------------------------------
# Edit summary
#   1) The table on slide 6 switches from the deck's custom table style to
#      the built-in table style {AD28F7E2-D3BD-4302-AB9C-D42BEB97D755}.
#   2) The presentation's theme colours change from the custom "Integral"
#      palette over to the standard Office palette (dk1/lt1 are already
#      black/white in both palettes, so only the other 10 slots move).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Table style -------------------------------------------------------
$tableSlide = $p.Slides.Item(6)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{AD28F7E2-D3BD-4302-AB9C-D42BEB97D755}", $false)

# --- 2) Theme colours ------------------------------------------------------
# Note: ThemeColor.RGB uses the OLE COLORREF byte order (0x00BBGGRR), i.e.
# the bytes of a plain RRGGBB hex colour are reversed before being assigned.
$themeColors = $p.Slides.Item(1).ThemeColorScheme

$themeColors.Colors(3).RGB  = 0x6A5444   # dk2      44546A
$themeColors.Colors(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$themeColors.Colors(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$themeColors.Colors(6).RGB  = 0x317DED   # accent2  ED7D31
$themeColors.Colors(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$themeColors.Colors(8).RGB  = 0x00C0FF   # accent4  FFC000
$themeColors.Colors(9).RGB  = 0xC47244   # accent5  4472C4
$themeColors.Colors(10).RGB = 0x47AD70   # accent6  70AD47
$themeColors.Colors(11).RGB = 0xC16305   # hlink    0563C1
$themeColors.Colors(12).RGB = 0x724F95   # folHlink 954F72
